# Automatische test-sync: 2025-06-19 22:30:50
#
# Appends the 22:30:10 "Factuur verzoek" mail-log entry as row 43 on the
# "Logs" sheet and refreshes the "Dashboard" pivot-style summary table to
# match: Klacht/Probleem and Offerte/Prijsaanvraag trade places (both sit
# at 4) and Factuur/Administratie climbs from 3 to 4.

$wb = $excel.ActiveWorkbook

# --- Logs: append the new mail entry on row 43 -----------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(43, 1).Value = "Factuur verzoek"
$logs.Cells.Item(43, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(43, 3).Value = "Kunt u mij de factuur van mijn laatste bestelling toesturen?"
$logs.Cells.Item(43, 4).Value = "Factuur / Administratie"
# column E (Antwoord) intentionally left blank - no reply yet
$logs.Cells.Item(43, 6).Value = "2025-06-19 22:30:10"
$logs.Cells.Item(43, 7).Value = "Nee"

# --- Dashboard: swap the row 5/6 category labels, bump the Factuur count ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(6, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(7, 2).Value = 4
